# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (column E) and
# "Correspond Handback DateTime" (column H) values for the 67b25730... row
# (row 4) on both the "zh-cn" and "de-de" sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E4").Value = "2016-03-11 08:26:37"
$wsZhCn.Range("H4").Value = "2016-03-11 08:26:54"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E4").Value = "2016-03-11 08:26:40"
$wsDeDe.Range("H4").Value = "2016-03-11 08:26:59"
